$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header row grows taller (wrapped footnote reference likely added to the title)
$ws.Rows.Item(1).RowHeight = 39.2

# The "female" count for the (Missing) species row now carries a footnote marker (a)
$ws.Range("B4").Value = "female (a)"

# Give the Total row (27) the same explicit row height as the rest of the table
$ws.Rows.Item(27).RowHeight = 27

# Add a new footnote row (28) below the Total row, matching its formatting
$ws.Range("A27:F27").Copy($ws.Range("A28:F28"))
$ws.Range("A28").Value = "(a) Some comment."
$ws.Range("B28:F28").ClearContents()
$ws.Rows.Item(28).RowHeight = 27
